$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2472.625
$ws.Range("I132").Value = 1991.3462
$ws.Range("J132").Value = 3366.4285
$ws.Range("K132").Value = 5974.0386
$ws.Range("L132").Value = 10099.2855
$ws.Range("M132").Value = -3444.0386
$ws.Range("N132").Value = -15159.2855

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1333.5333
$ws.Range("I2").Value = 829
$ws.Range("J2").Value = 2342.6
$ws.Range("K2").Value = 829
$ws.Range("L2").Value = 2342.6
$ws.Range("M2").Value = -716
$ws.Range("N2").Value = -2568.6
$ws.Range("H32").Value = 23921.629
$ws.Range("I32").Value = 25575.232
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 25575.232
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -25288.232
$ws.Range("N32").Value = -14574
$ws.Range("H109").Value = 32377
$ws.Range("J109").Value = 32377
$ws.Range("L109").Value = 32377
$ws.Range("N109").Value = -35151
$ws.Range("H116").Value = 1333.5333
$ws.Range("I116").Value = 829
$ws.Range("J116").Value = 2342.6
$ws.Range("K116").Value = 829
$ws.Range("L116").Value = 2342.6
$ws.Range("M116").Value = 1465
$ws.Range("N116").Value = -6930.6
$ws.Range("H132").Value = 1650.8823
$ws.Range("I132").Value = 1441.625
$ws.Range("K132").Value = 4324.875
$ws.Range("M132").Value = -1794.875
$ws.Range("H139").Value = 71275
$ws.Range("J139").Value = 71275
$ws.Range("L139").Value = 71275
$ws.Range("N139").Value = -81555

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1333.5333
$ws.Range("I3").Value = 829
$ws.Range("J3").Value = 2342.6
$ws.Range("K3").Value = 829
$ws.Range("L3").Value = 2342.6
$ws.Range("M3").Value = -715
$ws.Range("N3").Value = -2570.6
$ws.Range("H7").Value = 45751.5
$ws.Range("I7").Value = 43003
$ws.Range("K7").Value = 43003
$ws.Range("M7").Value = -42890
$ws.Range("H100").Value = 23731.5
$ws.Range("J100").Value = 23731.5
$ws.Range("L100").Value = 23731.5
$ws.Range("N100").Value = -25895.5
$ws.Range("H108").Value = 50101
$ws.Range("J108").Value = 50101
$ws.Range("L108").Value = 50101
$ws.Range("N108").Value = -57781
$ws.Range("H134").Value = 2643.4546
$ws.Range("I134").Value = 2429.5789
$ws.Range("J134").Value = 3998
$ws.Range("K134").Value = 7288.736699999999
$ws.Range("L134").Value = 11994
$ws.Range("M134").Value = -4753.736699999999
$ws.Range("N134").Value = -17064

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38465590
$ws.Range("I31").Value = 58826892
$ws.Range("K31").Value = 58826892
$ws.Range("M31").Value = -58826597
$ws.Range("H34").Value = 38465590
$ws.Range("I34").Value = 58826892
$ws.Range("K34").Value = 58826892
$ws.Range("M34").Value = -58826690
$ws.Range("H137").Value = 49081.11
$ws.Range("J137").Value = 64346
$ws.Range("L137").Value = 64346
$ws.Range("N137").Value = -74546
$ws.Range("H138").Value = 39941.25
$ws.Range("J138").Value = 39941.25
$ws.Range("L138").Value = 39941.25
$ws.Range("N138").Value = -50221.25
$ws.Range("H140").Value = 78659.28999999999
$ws.Range("J140").Value = 78659.28999999999
$ws.Range("L140").Value = 78659.28999999999
$ws.Range("N140").Value = -89019.28999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2465
$ws.Range("I46").Value = 358.7143
$ws.Range("J46").Value = 4571.2856
$ws.Range("K46").Value = 1076.1429
$ws.Range("L46").Value = 13713.8568
$ws.Range("M46").Value = -985.1428999999998
$ws.Range("N46").Value = -13895.8568
$ws.Range("H68").Value = 1257.826
$ws.Range("J68").Value = 1231.4117
$ws.Range("L68").Value = 3694.2351
$ws.Range("N68").Value = -5316.2351
$ws.Range("H71").Value = 1257.826
$ws.Range("J71").Value = 1231.4117
$ws.Range("L71").Value = 11082.7053
$ws.Range("N71").Value = -19194.7053
$ws.Range("H113").Value = 659.6667
$ws.Range("I113").Value = 501.5
$ws.Range("J113").Value = 712.3889
$ws.Range("K113").Value = 1504.5
$ws.Range("L113").Value = 2137.1667
$ws.Range("M113").Value = 665.5
$ws.Range("N113").Value = -6477.1667
$ws.Range("H121").Value = 56657
$ws.Range("J121").Value = 72390.42999999999
$ws.Range("L121").Value = 217171.29
$ws.Range("N121").Value = -219791.29
$ws.Range("H132").Value = 1735.4736
$ws.Range("J132").Value = 2036.8462
$ws.Range("L132").Value = 18331.6158
$ws.Range("N132").Value = -23391.6158

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 23769.23
$ws.Range("I5").Value = 9750
$ws.Range("K5").Value = 9750
$ws.Range("M5").Value = -9638
$ws.Range("H138").Value = 29995
$ws.Range("J138").Value = 29995
$ws.Range("L138").Value = 29995
$ws.Range("N138").Value = -40275

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9606.666999999999
$ws.Range("J2").Value = 9962
$ws.Range("L2").Value = 9962
$ws.Range("N2").Value = -10186
$ws.Range("H36").Value = 45000
$ws.Range("J36").Value = 45000
$ws.Range("L36").Value = 45000
$ws.Range("N36").Value = -46124
$ws.Range("H98").Value = 33500
$ws.Range("J98").Value = 33500
$ws.Range("L98").Value = 33500
$ws.Range("N98").Value = -39490
$ws.Range("H132").Value = 4888.914
$ws.Range("I132").Value = 4754.125
$ws.Range("K132").Value = 14262.375
$ws.Range("M132").Value = -11732.375
$ws.Range("H136").Value = 3799.8372
$ws.Range("I136").Value = 4177.5835
$ws.Range("K136").Value = 12532.7505
$ws.Range("M136").Value = -9982.750499999998
$ws.Range("H137").Value = 83900
$ws.Range("J137").Value = 83900
$ws.Range("L137").Value = 83900
$ws.Range("N137").Value = -94100
$ws.Range("H139").Value = 40056.332
$ws.Range("J139").Value = 40056.332
$ws.Range("L139").Value = 40056.332
$ws.Range("N139").Value = -50336.332

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2051.1904
$ws.Range("I132").Value = 2594.4
$ws.Range("J132").Value = 693.1667
$ws.Range("K132").Value = 7783.200000000001
$ws.Range("L132").Value = 2079.5001
$ws.Range("M132").Value = -5253.200000000001
$ws.Range("N132").Value = -7139.5001
